$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value = "31/12/2025"
$ws.Range("C7").Value = 76.81774852108779

$ws.Range("B13").Value = "31/12/2025"
$ws.Range("C13").Value = 76.3785534556983

$ws.Range("B19").Value = "31/12/2025"
$ws.Range("C19").Value = 76.57679819713221
